$wb = $excel.ActiveWorkbook

# Rename sheets: "Repeaters" -> "RepeatersOld", "Repeaters Updated" -> "Repeaters"
$wsOld = $wb.Worksheets.Item("Repeaters")
$wsOld.Name = "RepeatersOld"

$wsNew = $wb.Worksheets.Item("Repeaters Updated")
$wsNew.Name = "Repeaters"

# Make the (renamed) "Repeaters" sheet the active tab, with A15 selected
$wsNew.Activate()
$wsNew.Range("A15").Select()
